$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "Remédios "
$ws.Range("A21").Value = "Mercadinho condomínio "
$ws.Range("A20:A21").Style = "Normal"
